$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 824875.5
$ws.Range("J17").Value = 824875.5
$ws.Range("L17").Value = 2474626.5
$ws.Range("N17").Value = -2474962.5
$ws.Range("H92").Value = 684.2
$ws.Range("I92").Value = 606.8889
$ws.Range("J92").Value = 800.1667
$ws.Range("K92").Value = 606.8889
$ws.Range("L92").Value = 800.1667
$ws.Range("M92").Value = 641.1111
$ws.Range("N92").Value = -3296.1667
$ws.Range("H101").Value = 829.9091
$ws.Range("I101").Value = 638.2222
$ws.Range("J101").Value = 1692.5
$ws.Range("K101").Value = 1914.6666
$ws.Range("L101").Value = 5077.5
$ws.Range("M101").Value = -292.6666
$ws.Range("N101").Value = -8321.5
$ws.Range("H112").Value = 5496540.5
$ws.Range("J112").Value = 5496540.5
$ws.Range("L112").Value = 16489621.5
$ws.Range("N112").Value = -16491837.5
$ws.Range("H135").Value = 13159035
$ws.Range("I135").Value = 851.1539
$ws.Range("K135").Value = 7660.3851
$ws.Range("M135").Value = -5125.3851
$ws.Range("H138").Value = 9010703
$ws.Range("I138").Value = 1132.2222
$ws.Range("J138").Value = 17546086
$ws.Range("K138").Value = 3396.6666
$ws.Range("L138").Value = 52638258
$ws.Range("M138").Value = 1743.3334
$ws.Range("N138").Value = -52648538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1177.3334
$ws.Range("I2").Value = 1021.5
$ws.Range("K2").Value = 1021.5
$ws.Range("M2").Value = -908.5
$ws.Range("H32").Value = 12992463
$ws.Range("I32").Value = 15387424
$ws.Range("K32").Value = 15387424
$ws.Range("M32").Value = -15387137
$ws.Range("H61").Value = 24393822
$ws.Range("I61").Value = 31252036
$ws.Range("J61").Value = 9056.777
$ws.Range("K61").Value = 31252036
$ws.Range("L61").Value = 9056.777
$ws.Range("M61").Value = -31251824
$ws.Range("N61").Value = -9480.777
$ws.Range("H97").Value = 1606.4783
$ws.Range("I97").Value = 1577
$ws.Range("K97").Value = 1577
$ws.Range("M97").Value = -1081
$ws.Range("H102").Value = 2644.4546
$ws.Range("I102").Value = 1660.7059
$ws.Range("K102").Value = 1660.7059
$ws.Range("M102").Value = -38.70589999999993
$ws.Range("H116").Value = 1177.3334
$ws.Range("I116").Value = 1021.5
$ws.Range("K116").Value = 1021.5
$ws.Range("M116").Value = 1272.5
$ws.Range("H132").Value = 33336782
$ws.Range("I132").Value = 3572
$ws.Range("J132").Value = 333335680
$ws.Range("K132").Value = 10716
$ws.Range("L132").Value = 1000007040
$ws.Range("M132").Value = -8186
$ws.Range("N132").Value = -1000012100
$ws.Range("H136").Value = 24393822
$ws.Range("I136").Value = 31252036
$ws.Range("J136").Value = 9056.777
$ws.Range("K136").Value = 93756108
$ws.Range("L136").Value = 27170.331
$ws.Range("M136").Value = -93753558
$ws.Range("N136").Value = -32270.331
$ws.Range("H139").Value = 115715
$ws.Range("J139").Value = 115715
$ws.Range("L139").Value = 115715
$ws.Range("N139").Value = -125995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1177.3334
$ws.Range("I3").Value = 1021.5
$ws.Range("K3").Value = 1021.5
$ws.Range("M3").Value = -907.5
$ws.Range("H99").Value = 5654.1113
$ws.Range("I99").Value = 3975.2
$ws.Range("K99").Value = 3975.2
$ws.Range("M99").Value = -2477.2
$ws.Range("H105").Value = 10995.8
$ws.Range("I105").Value = 13419.75
$ws.Range("K105").Value = 13419.75
$ws.Range("M105").Value = -11672.75
$ws.Range("H134").Value = 2542.451
$ws.Range("I134").Value = 2322.7021
$ws.Range("J134").Value = 5124.5
$ws.Range("K134").Value = 6968.106299999999
$ws.Range("L134").Value = 15373.5
$ws.Range("M134").Value = -4433.106299999999
$ws.Range("N134").Value = -20443.5
$ws.Range("H137").Value = 195000
$ws.Range("J137").Value = 195000
$ws.Range("L137").Value = 195000
$ws.Range("N137").Value = -205200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 10402.9
$ws.Range("I22").Value = 16897
$ws.Range("J22").Value = 661.75
$ws.Range("K22").Value = 16897
$ws.Range("L22").Value = 661.75
$ws.Range("M22").Value = -16547
$ws.Range("N22").Value = -1361.75
$ws.Range("H31").Value = 18185426
$ws.Range("I31").Value = 2458.2856
$ws.Range("J31").Value = 166679660
$ws.Range("K31").Value = 2458.2856
$ws.Range("L31").Value = 166679660
$ws.Range("M31").Value = -2163.2856
$ws.Range("N31").Value = -166680250
$ws.Range("H34").Value = 18185426
$ws.Range("I34").Value = 2458.2856
$ws.Range("J34").Value = 166679660
$ws.Range("K34").Value = 2458.2856
$ws.Range("L34").Value = 166679660
$ws.Range("M34").Value = -2256.2856
$ws.Range("N34").Value = -166680064
$ws.Range("H58").Value = 2190.524
$ws.Range("J58").Value = 4500
$ws.Range("L58").Value = 4500
$ws.Range("N58").Value = -4906
$ws.Range("H136").Value = 2190.524
$ws.Range("J136").Value = 4500
$ws.Range("L136").Value = 13500
$ws.Range("N136").Value = -18600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 198999
$ws.Range("J37").Value = 198999
$ws.Range("L37").Value = 596997
$ws.Range("N37").Value = -597221
$ws.Range("H56").Value = 23767.2
$ws.Range("I56").Value = 23767.2
$ws.Range("K56").Value = 23767.2
$ws.Range("M56").Value = -23237.2
$ws.Range("H99").Value = 1178.6666
$ws.Range("I99").Value = 1178.6666
$ws.Range("K99").Value = 3535.9998
$ws.Range("M99").Value = -1289.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6308.625
$ws.Range("I70").Value = 4379
$ws.Range("J70").Value = 8238.25
$ws.Range("K70").Value = 4379
$ws.Range("L70").Value = 8238.25
$ws.Range("M70").Value = -4109
$ws.Range("N70").Value = -8778.25
$ws.Range("H73").Value = 6308.625
$ws.Range("I73").Value = 4379
$ws.Range("J73").Value = 8238.25
$ws.Range("K73").Value = 4379
$ws.Range("L73").Value = 8238.25
$ws.Range("M73").Value = -3443
$ws.Range("N73").Value = -10110.25
$ws.Range("H102").Value = 4462.7334
$ws.Range("I102").Value = 4104.5557
$ws.Range("K102").Value = 4104.5557
$ws.Range("M102").Value = -2482.5557
$ws.Range("H122").Value = 6818.25
$ws.Range("J122").Value = 9719
$ws.Range("L122").Value = 29157
$ws.Range("N122").Value = -34057
$ws.Range("H126").Value = 27280294
$ws.Range("I126").Value = 25010960
$ws.Range("K126").Value = 75032880
$ws.Range("M126").Value = -75030410
$ws.Range("H132").Value = 3181.3462
$ws.Range("I132").Value = 2986.1428
$ws.Range("K132").Value = 8958.428400000001
$ws.Range("M132").Value = -6428.428400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2645.5
$ws.Range("I22").Value = 654.6
$ws.Range("J22").Value = 3411.2307
$ws.Range("K22").Value = 654.6
$ws.Range("L22").Value = 3411.2307
$ws.Range("M22").Value = -359.6
$ws.Range("N22").Value = -4001.2307
$ws.Range("H27").Value = 2645.5
$ws.Range("I27").Value = 654.6
$ws.Range("J27").Value = 3411.2307
$ws.Range("K27").Value = 654.6
$ws.Range("L27").Value = 3411.2307
$ws.Range("M27").Value = -547.6
$ws.Range("N27").Value = -3625.2307
$ws.Range("H46").Value = 1122.5942
$ws.Range("I46").Value = 611.5472
$ws.Range("K46").Value = 611.5472
$ws.Range("M46").Value = -423.5472
$ws.Range("H68").Value = 3093.25
$ws.Range("I68").Value = 2835.4443
$ws.Range("K68").Value = 2835.4443
$ws.Range("M68").Value = -2086.4443
$ws.Range("H71").Value = 3093.25
$ws.Range("I71").Value = 2835.4443
$ws.Range("K71").Value = 14177.2215
$ws.Range("M71").Value = -10433.2215
$ws.Range("H122").Value = 4505.6924
$ws.Range("I122").Value = 4127.136
$ws.Range("J122").Value = 4995.5884
$ws.Range("K122").Value = 12381.408
$ws.Range("L122").Value = 14986.7652
$ws.Range("M122").Value = -9931.408000000001
$ws.Range("N122").Value = -19886.7652
$ws.Range("H132").Value = 86959390
$ws.Range("I132").Value = 2809.3572
$ws.Range("K132").Value = 8428.071599999999
$ws.Range("M132").Value = -5898.071599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 667959.8
$ws.Range("I100").Value = 715528.4
$ws.Range("K100").Value = 1431056.8
$ws.Range("M100").Value = -1430515.8
$ws.Range("H113").Value = 752.9524
$ws.Range("I113").Value = 305.16666
$ws.Range("J113").Value = 1350
$ws.Range("K113").Value = 915.4999799999999
$ws.Range("L113").Value = 4050
$ws.Range("M113").Value = 1254.50002
$ws.Range("N113").Value = -8390
